# Upload new version with timestamp
# - Removes the "ترمومتر ديجتال" (Digital Thermometer) row
# - Removes the "سرنجات انسولين" (Insulin syringes) row
# - Updates the running total accordingly
# - Bumps the generated-at timestamp from 8:54 PM to 8:55 PM

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "ترمومتر ديجتال" product row (row 43).
$ws.Range("A43:Q43").EntireRow.Delete()

# After the row above shifts everything up by one, the
# "سرنجات انسولين" product row is now at row 50.
$ws.Range("A50:Q50").EntireRow.Delete()

# Update the grand total (two rows, now at row 52, shifted up from 54)
# to reflect the removal of the two rows (50.00 + 14.00 = 64.00 less).
$ws.Range("P52").Value = 2050.3049999999998

# Update the generated timestamp shown in the footer (now at row 53).
$ws.Range("A53").Value = "Saturday, 7 June, 2025 8:55 PM"
